$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Priority column (E) for rows 4-7 moves from "low" to "ht" on both language sheets.
$wsZhCn.Range("E4:E7").Value = "ht"
$wsDeDe.Range("E4:E7").Value = "ht"

# Latest Handoff Datetime column (H) for rows 4-7 is refreshed with a new
# handoff report generation timestamp on both language sheets.
$wsZhCn.Range("H4:H7").Value = "2016-08-31 08:37:36"
$wsDeDe.Range("H4:H7").Value = "2016-08-31 08:37:42"

# The Overview sheet mirrors the latest handoff generation date (the de-de
# timestamp), so refresh it to match the new handoff run too.
$wsOverview.Range("G4:G7").Value = "2016-08-31 08:37:42"
